$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price / volume(1h) values scraped on Fri Jun 23 23:16:40 UTC 2023
$ws.Range("D2").Formula = '''30.656.37'
$ws.Range("E2").Formula = '''  +2.35%  '
$ws.Range("D3").Formula = '''1.892.64'
$ws.Range("E3").Formula = '''  +1.05%  '
$ws.Range("E4").Formula = '''  +0.10%  '
$ws.Range("D5").Formula = '''244.48'
$ws.Range("E5").Formula = '''  +1.74%  '
$ws.Range("E6").Formula = '''  +0.04%  '
$ws.Range("D7").Formula = '''0.4971'
$ws.Range("E7").Formula = '''  +0.43%  '
$ws.Range("D8").Formula = '''0.2963'
$ws.Range("E8").Formula = '''  +2.16%  '
$ws.Range("D9").Formula = '''0.06817'
$ws.Range("E9").Formula = '''  +3.92%  '
$ws.Range("D10").Formula = '''1.891.94'
$ws.Range("E10").Formula = '''  +1.03%  '
$ws.Range("D11").Formula = '''17.07'
$ws.Range("E11").Formula = '''  +2.52%  '
$ws.Range("D12").Formula = '''0.07328'
$ws.Range("E12").Formula = '''  +2.24%  '
$ws.Range("D13").Formula = '''91.41'
$ws.Range("E13").Formula = '''  +6.61%  '
$ws.Range("D14").Formula = '''5.090'
$ws.Range("E14").Formula = '''  +5.28%  '
$ws.Range("D15").Formula = '''0.6745'
$ws.Range("E15").Formula = '''  +2.86%  '
$ws.Range("D16").Formula = '''30.654.86'
$ws.Range("E16").Formula = '''  +2.43%  '
$ws.Range("D17").Formula = '''0.000007935'
$ws.Range("E17").Formula = '''  +1.43%  '
$ws.Range("E18").Formula = '''  -0.03%  '
$ws.Range("D19").Formula = '''13.27'
$ws.Range("E19").Formula = '''  +4.96%  '
$ws.Range("D20").Formula = '''2.134.69'
$ws.Range("E20").Formula = '''  +0.83%  '
$ws.Range("D21").Formula = '''1.000'
$ws.Range("E21").Formula = '''  +0.02%  '
$ws.Range("D22").Formula = '''4.868'
$ws.Range("E22").Formula = '''  +3.10%  '
$ws.Range("D23").Formula = '''177.34'
$ws.Range("E23").Formula = '''  +31.91%  '
$ws.Range("D24").Formula = '''6.065'
$ws.Range("E24").Formula = '''  +9.11%  '
$ws.Range("D25").Formula = '''9.300'
$ws.Range("E25").Formula = '''  +3.21%  '
$ws.Range("D26").Formula = '''154.60'
$ws.Range("E26").Formula = '''  +3.26%  '
$ws.Range("D27").Formula = '''18.79'
$ws.Range("E27").Formula = '''  +13.05%  '
$ws.Range("D28").Formula = '''1.928'
$ws.Range("E28").Formula = '''  +2.09%  '
$ws.Range("D29").Formula = '''1.388'
$ws.Range("E29").Formula = '''  +1.14%  '
$ws.Range("D30").Formula = '''4.344'
$ws.Range("E30").Formula = '''  +4.98%  '
$ws.Range("D31").Formula = '''0.08948'
$ws.Range("E31").Formula = '''  +3.30%  '
$ws.Range("D32").Formula = '''4.038'
$ws.Range("E32").Formula = '''  +3.16%  '
$ws.Range("D33").Formula = '''0.05200'
$ws.Range("E33").Formula = '''  +3.98%  '
$ws.Range("D34").Formula = '''0.7398'
$ws.Range("E34").Formula = '''  +6.45%  '
$ws.Range("D35").Formula = '''1.136'
$ws.Range("E35").Formula = '''  +4.55%  '
$ws.Range("E36").Formula = '''  +0.77%  '
$ws.Range("D37").Formula = '''0.01876'
$ws.Range("E37").Formula = '''  +10.85%  '
$ws.Range("D38").Formula = '''2.702'
$ws.Range("E38").Formula = '''  +0.76%  '
$ws.Range("D39").Formula = '''2.173'
$ws.Range("E39").Formula = '''  +0.78%  '
$ws.Range("D40").Formula = '''0.9332'
$ws.Range("E40").Formula = '''  +1.14%  '
$ws.Range("D41").Formula = '''0.4360'
$ws.Range("E41").Formula = '''  +4.97%  '
$ws.Range("D42").Formula = '''106.05'
$ws.Range("E42").Formula = '''  +4.61%  '
$ws.Range("D43").Formula = '''5.806'
$ws.Range("E43").Formula = '''  -2.22%  '
$ws.Range("E44").Formula = '''  +0.20%  '
$ws.Range("D45").Formula = '''7.671'
$ws.Range("E45").Formula = '''  +4.51%  '
$ws.Range("D46").Formula = '''0.1354'
$ws.Range("E46").Formula = '''  +8.14%  '
$ws.Range("D47").Formula = '''0.05847'
$ws.Range("E47").Formula = '''  +3.51%  '
$ws.Range("D48").Formula = '''33.43'
$ws.Range("E48").Formula = '''  +3.49%  '
$ws.Range("D49").Formula = '''0.3895'
$ws.Range("E49").Formula = '''  +6.05%  '
$ws.Range("D50").Formula = '''8.482'
$ws.Range("E50").Formula = '''  +5.45%  '
$ws.Range("D51").Formula = '''1.384'
$ws.Range("E51").Formula = '''  +4.30%  '
